# Fix OCR output text in the "metadata" sheet (Key/Value table).
# These corrections patch a handful of misread characters coming out of
# the CIMB statement OCR pass (account number, account holder name and
# address) before the data is exported to the Eagle Eye Excel format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# Each OCR'd value ends with a form-feed control character (rendered in
# the OOXML as the literal escape "_x000C_") following a newline - keep
# that trailer intact while only touching the mis-OCR'd characters.
$ff = [char]12

# B2 - account_number: "T6-1806128-3" -> "T6-18060128-3"
$ws.Cells.Item(2, 2).Value = "T6-18060128-3`n" + $ff

# B3 - account_holder: "NUR IZZAHTI BINTI AZEMAN" -> "NUR IZZAHTI BINT! AZEMAN"
$ws.Cells.Item(3, 2).Value = "NUR IZZAHTI BINT! AZEMAN`n" + $ff

# B4 - address: "53100 SELANGOR SELANGOR" -> "53100 SELANGOR, SELANGOR"
$ws.Cells.Item(4, 2).Value = "NO 19 JALAN Nd FASA DA`nTAMAN MELAWATI`n`nKUALA LUMPUR`n`n53100 SELANGOR, SELANGOR`n" + $ff

# B5 - account_holder: "3419 JALAN Ald FASA OA" -> "7419 JALAN Ald FASA GA"
$ws.Cells.Item(5, 2).Value = "UR IZZAHTI BINT! AZEM`n7419 JALAN Ald FASA GA`n" + $ff
